$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the whole table down by one row (new row 1 inserted above the old header).
$ws.Rows("1:1").Insert()

# --- New "Mögliche Lösung" (solution) column G ---------------------------------

# B7: new bug row "Infinity-Jump-Bug"
$ws.Range("B7").Value = "Infinity-Jump-Bug"

# G2: new column header
$ws.Range("G2").Value = "Mögliche Lösung"

# C5: replace old "Kann gelöst..." text with a short dash (that text is reused below in G5)
$ws.Range("C5").Value = "-"

# G4: solution text for row 4 (Dobble Jump duplicate GameObject bug)
$ws.Range("G4").Value = "das GameObject Inventory_UI und Inventory_Main_Driver werden 2x geladen checken ob es schon geladen ist"

# G3: solution text for row 3 (Musikbug)
$ws.Range("G3").Value = "Checken ob GameObject schon erstellt wurde"

# C7: description for the new Infinity-Jump-Bug row
$ws.Range("C7").Value = "Entsteht wenn man bei einem Hinweis"

# B8: new bug row "Bugged-Knockback"
$ws.Range("B8").Value = "Bugged-Knockback"

# C8: description for Bugged-Knockback
$ws.Range("C8").Value = "Entsteht wenn Player schaden nimmt"

# D8: cause for Bugged-Knockback
$ws.Range("D8").Value = "Wenn der Player schaden erleidet"

# G5 reuses the text originally in C4 ("Kann gelöst werden indehm der Player mid Tag gesucht wird")
$ws.Range("G5").Value = "Kann gelöst werden indehm der Player mid Tag gesucht wird"

# --- Fill in the "Solved" (E) column 0s that are now part of populated rows ----
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("E8").Value = 0

# --- Formatting: give column G the same look as the existing bordered columns --
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

$ws.Range("F3").Copy()
$ws.Range("G3:G5").PasteSpecial(-4122)

$ws.Range("F6:F32").Copy()
$ws.Range("G6:G32").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Column G width to roughly match the longest "Mögliche Lösung" entry.
$ws.Columns("G").ColumnWidth = 100.28515625

# Park the selection back on A1 (closest match to the saved view having no explicit selection).
$ws.Range("A1").Select() | Out-Null
